{"js": "// Insert \"Approximate \" as a new, separately-formatted run immediately\n// before the existing \"Algorithmic Image Matching to Reduce Online\" run\n// in the thesis title, so the title reads:\n//   \"Approximate Algorithmic Image Matching to Reduce Online\n//    Storage Overhead of User Submitted Images\"\n// The new run is given the same bold / italic / 16pt (sz=32 half-points)\n// character formatting as the run it precedes.\n\nconst results = context.document.body.search(\n  \"Algorithmic Image Matching to Reduce Online\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find title text \"Algorithmic Image Matching to Reduce Online\".');\n}\n\nconst titleRun = results.items[0];\n// Get a zero-length range sitting right before the matched text so the new\n// text lands in its own run instead of being merged into titleRun's run.\nconst insertionPoint = titleRun.getRange(Word.RangeLocation.start);\nconst inserted = insertionPoint.insertText(\"Approximate \", Word.InsertLocation.before);\n\n// Match the title run's character formatting (bold, italic, 32 half-points = 16pt).\ninserted.font.bold = true;\ninserted.font.italic = true;\ninserted.font.size = 16;\n\nawait context.sync();\n", "ps1": "# Insert \"Approximate \" as a new, separately-formatted run immediately\n# before the existing \"Algorithmic Image Matching to Reduce Online\" run\n# in the thesis title, so the title reads:\n#   \"Approximate Algorithmic Image Matching to Reduce Online\n#    Storage Overhead of User Submitted Images\"\n# The new run is given the same bold / italic / 16pt (sz=32 half-points)\n# character formatting as the run it precedes.\n\n$d = $word.ActiveDocument\n\n$titleRange = $d.Content\n$found = $titleRange.Find.Execute(\"Algorithmic Image Matching to Reduce Online\")\nif (-not $found) {\n    throw 'Could not find title text \"Algorithmic Image Matching to Reduce Online\".'\n}\n\n# Collapse to the start of the match so the new text is inserted as its own\n# run instead of being merged into the matched run.\n$titleRange.Collapse(1)  # wdCollapseStart\n$titleRange.InsertBefore(\"Approximate \")\n\n# Re-find the newly inserted run and give it the same character formatting\n# (bold, italic, 16pt) as the title text that follows it.\n$newRunRange = $d.Content\n$newRunRange.Find.Execute(\"Approximate \")\n$newRunRange.Font.Bold = 1\n$newRunRange.Font.Italic = 1\n$newRunRange.Font.Size = 16\n"}
